{"js": "// Accept all tracked changes (insertions/deletions/formatting changes) in the document.\nconst trackedChanges = context.document.getTrackedChanges();\ntrackedChanges.load(\"items\");\nawait context.sync();\n\ntrackedChanges.acceptAll();\nawait context.sync();\n", "ps1": "# Accept all tracked changes (insertions/deletions/formatting changes) in the document.\n$d = $word.ActiveDocument\n$d.Revisions.AcceptAll()\n"}
